$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: update styles (A11 -> style used by rows 12-14 "A" column, B11 -> style used by rows 12-14 "B" column) ---
# Existing style for A11 was the "numFmt164 center/center" style; new style s="7" is
# "numFmt164 vertical-center only" (no horizontal center). We replicate via range.Style
# using a newly defined named style mapped to the same formatting as the new cellXfs entry.

# Create (or reuse) a style that matches cellXfs index 7: numFmtId 164, font index 1 (Calibri Light),
# vertical center alignment only (no horizontal alignment).
$styleName = "CalendarVCenter"
$styleExists = $false
foreach ($s in $wb.Styles) {
    if ($s.Name -eq $styleName) { $styleExists = $true }
}
if (-not $styleExists) {
    $newStyle = $wb.Styles.Add($styleName)
    $newStyle.VerticalAlignment = -4108  # xlVAlignCenter
    $newStyle.HorizontalAlignment = 1    # xlHAlignGeneral
    $newStyle.NumberFormat = "[$-3000401]0"
    $newStyle.Font.Name = "Calibri Light"
}

$ws.Range("A11").Style = $styleName
$ws.Range("A12").Style = $styleName
$ws.Range("A13").Style = $styleName

# B11:B14 uses style s="5" which matches the existing style of B2/B5/B8 (fontId 1, center/center)
$bStyleName = "CalendarCenter"
$bStyleExists = $false
foreach ($s in $wb.Styles) {
    if ($s.Name -eq $bStyleName) { $bStyleExists = $true }
}
if (-not $bStyleExists) {
    $newBStyle = $wb.Styles.Add($bStyleName)
    $newBStyle.VerticalAlignment = -4108  # xlVAlignCenter
    $newBStyle.HorizontalAlignment = -4108 # xlHAlignCenter
    $newBStyle.Font.Name = "Calibri Light"
}

$ws.Range("B11:B14").Style = $bStyleName

# --- New rows 12-14 content ---
$ws.Range("C12").Value = "بهتر کردن شکل ظاهری تقویم و تبدیل table به div"
$ws.Range("D12").Value = "1403/09/15"
$ws.Range("E12").Value = 14
$ws.Range("F12").Value = 16

$ws.Range("C13").Value = "اضافه کردن دو سال دیگر و تغییر عملکرد بر اساس آن ها و بهبود ظاهر"
$ws.Range("D13").Value = "1403/09/16"
$ws.Range("E13").Value = 14.5
$ws.Range("F13").Value = 17.5

$ws.Range("C14").Value = "اضافه کردن فونت فارسی مناسب "
$ws.Range("D14").Value = "1403/09/18"
$ws.Range("E14").Value = 17.5
$ws.Range("F14").Value = 21.5

# Apply style s="1" (center/center, fontId 1) to C/D/E/F of rows 12-14, matching other data rows
$cdefStyleName = "CalendarCDEF"
$cdefExists = $false
foreach ($s in $wb.Styles) {
    if ($s.Name -eq $cdefStyleName) { $cdefExists = $true }
}
if (-not $cdefExists) {
    $newCdefStyle = $wb.Styles.Add($cdefStyleName)
    $newCdefStyle.VerticalAlignment = -4108
    $newCdefStyle.HorizontalAlignment = -4108
    $newCdefStyle.Font.Name = "Calibri Light"
}
$ws.Range("C12:F14").Style = $cdefStyleName

# --- Extend shared formula G2:G14 ---
$ws.Range("G12").Formula = "=F12-E12"
$ws.Range("G13").Formula = "=F13-E13"
$ws.Range("G14").Formula = "=F14-E14"
$ws.Range("G12:G14").Style = $cdefStyleName

# --- Merge B11:B14 ---
$ws.Range("B11:B14").Merge()

# --- Selection update ---
$ws.Range("F14").Select()
